$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.055.86'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +5.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.467.69'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +5.97%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '186.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +7.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '546.67'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.22%  '
$ws.Range('E7').Value = '  +2.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.465.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +6.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('E10').Value = '  +6.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.12'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('E12').Value = '  +11.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +7.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.41'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.019.51'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +5.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.466.58'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +5.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.607.51'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +6.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.121'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +4.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.25'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.78'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +7.62%  '
$ws.Range('E21').Value = '  +6.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '404.88'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +9.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.00'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +10.35%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.90'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.92%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.64'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.20'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +9.36%  '
$ws.Range('E27').Value = '  +10.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.26'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.77'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +4.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.62'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.14'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '680.26'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +6.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.87'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.64'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.66%  '
$ws.Range('E35').Value = '  +5.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.98'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0826'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +20.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '38.43'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +6.43%  '
$ws.Range('E39').Value = '  +5.45%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.40'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +24.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.79'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +15.21%  '
$ws.Range('E43').Value = '  +7.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.041.25'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.07%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.97'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +12.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0420'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +7.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.28'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.44%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.13'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +18.06%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.74'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.130'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.36%  '
